$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(105, 1).Value = 0
$ws.Cells.Item(105, 2).Value = 8
$ws.Cells.Item(105, 3).Value = 4.2935
$ws.Cells.Item(105, 4).Value = 0
$ws.Cells.Item(105, 5).Value = 250
$ws.Cells.Item(105, 6).Value = 112
$ws.Cells.Item(105, 7).Value = 0.16
$ws.Cells.Item(105, 8).Value = 0.158

$ws.Cells.Item(106, 1).Value = 1
$ws.Cells.Item(106, 2).Value = 8
$ws.Cells.Item(106, 3).Value = 1.7325
$ws.Cells.Item(106, 4).Value = 0
$ws.Cells.Item(106, 5).Value = 250
$ws.Cells.Item(106, 6).Value = 112
$ws.Cells.Item(106, 7).Value = 0.16
$ws.Cells.Item(106, 8).Value = 0.158

$ws.Cells.Item(107, 1).Value = -9
$ws.Cells.Item(107, 2).Value = 9
$ws.Cells.Item(107, 3).Value = 4.44
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 250
$ws.Cells.Item(107, 6).Value = 112
$ws.Cells.Item(107, 7).Value = 0.16
$ws.Cells.Item(107, 8).Value = 0.158

$ws.Cells.Item(108, 1).Value = -8
$ws.Cells.Item(108, 2).Value = 9
$ws.Cells.Item(108, 3).Value = 4.9847
$ws.Cells.Item(108, 4).Value = 0
$ws.Cells.Item(108, 5).Value = 250
$ws.Cells.Item(108, 6).Value = 112
$ws.Cells.Item(108, 7).Value = 0.16
$ws.Cells.Item(108, 8).Value = 0.158

$ws.Cells.Item(109, 1).Value = -7
$ws.Cells.Item(109, 2).Value = 9
$ws.Cells.Item(109, 3).Value = 5.4923
$ws.Cells.Item(109, 4).Value = 0
$ws.Cells.Item(109, 5).Value = 250
$ws.Cells.Item(109, 6).Value = 112
$ws.Cells.Item(109, 7).Value = 0.16
$ws.Cells.Item(109, 8).Value = 0.158

$ws.Cells.Item(110, 1).Value = -6
$ws.Cells.Item(110, 2).Value = 9
$ws.Cells.Item(110, 3).Value = 5.964
$ws.Cells.Item(110, 4).Value = 0
$ws.Cells.Item(110, 5).Value = 250
$ws.Cells.Item(110, 6).Value = 112
$ws.Cells.Item(110, 7).Value = 0.16
$ws.Cells.Item(110, 8).Value = 0.158

$ws.Cells.Item(111, 1).Value = -5
$ws.Cells.Item(111, 2).Value = 9
$ws.Cells.Item(111, 3).Value = 6.4409
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 250
$ws.Cells.Item(111, 6).Value = 112
$ws.Cells.Item(111, 7).Value = 0.16
$ws.Cells.Item(111, 8).Value = 0.158

$ws.Cells.Item(112, 1).Value = -4
$ws.Cells.Item(112, 2).Value = 9
$ws.Cells.Item(112, 3).Value = 6.6913
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 250
$ws.Cells.Item(112, 6).Value = 112
$ws.Cells.Item(112, 7).Value = 0.16
$ws.Cells.Item(112, 8).Value = 0.158

$ws.Cells.Item(113, 1).Value = -3
$ws.Cells.Item(113, 2).Value = 9
$ws.Cells.Item(113, 3).Value = 6.7729
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 250
$ws.Cells.Item(113, 6).Value = 112
$ws.Cells.Item(113, 7).Value = 0.16
$ws.Cells.Item(113, 8).Value = 0.158

$ws.Cells.Item(114, 1).Value = -2
$ws.Cells.Item(114, 2).Value = 9
$ws.Cells.Item(114, 3).Value = 6.5031
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 250
$ws.Cells.Item(114, 6).Value = 112
$ws.Cells.Item(114, 7).Value = 0.16
$ws.Cells.Item(114, 8).Value = 0.158

$ws.Cells.Item(115, 1).Value = -1
$ws.Cells.Item(115, 2).Value = 9
$ws.Cells.Item(115, 3).Value = 5.6733
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 250
$ws.Cells.Item(115, 6).Value = 112
$ws.Cells.Item(115, 7).Value = 0.16
$ws.Cells.Item(115, 8).Value = 0.158

$ws.Cells.Item(116, 1).Value = 0
$ws.Cells.Item(116, 2).Value = 9
$ws.Cells.Item(116, 3).Value = 4.6252
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 250
$ws.Cells.Item(116, 6).Value = 112
$ws.Cells.Item(116, 7).Value = 0.16
$ws.Cells.Item(116, 8).Value = 0.158

$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = 9
$ws.Cells.Item(117, 3).Value = 3.544
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 250
$ws.Cells.Item(117, 6).Value = 112
$ws.Cells.Item(117, 7).Value = 0.16
$ws.Cells.Item(117, 8).Value = 0.158

$ws.Cells.Item(118, 1).Value = 2
$ws.Cells.Item(118, 2).Value = 9
$ws.Cells.Item(118, 3).Value = 1.4599
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(118, 5).Value = 250
$ws.Cells.Item(118, 6).Value = 112
$ws.Cells.Item(118, 7).Value = 0.16
$ws.Cells.Item(118, 8).Value = 0.158

$ws.Cells.Item(119, 1).Value = -9
$ws.Cells.Item(119, 2).Value = 10
$ws.Cells.Item(119, 3).Value = 4.3695
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 250
$ws.Cells.Item(119, 6).Value = 112
$ws.Cells.Item(119, 7).Value = 0.16
$ws.Cells.Item(119, 8).Value = 0.158

$ws.Cells.Item(120, 1).Value = -8
$ws.Cells.Item(120, 2).Value = 10
$ws.Cells.Item(120, 3).Value = 4.9098
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 250
$ws.Cells.Item(120, 6).Value = 112
$ws.Cells.Item(120, 7).Value = 0.16
$ws.Cells.Item(120, 8).Value = 0.158

$ws.Cells.Item(121, 1).Value = -7
$ws.Cells.Item(121, 2).Value = 10
$ws.Cells.Item(121, 3).Value = 5.4203
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 250
$ws.Cells.Item(121, 6).Value = 112
$ws.Cells.Item(121, 7).Value = 0.16
$ws.Cells.Item(121, 8).Value = 0.158

$ws.Cells.Item(122, 1).Value = -6
$ws.Cells.Item(122, 2).Value = 10
$ws.Cells.Item(122, 3).Value = 5.9466
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 250
$ws.Cells.Item(122, 6).Value = 112
$ws.Cells.Item(122, 7).Value = 0.16
$ws.Cells.Item(122, 8).Value = 0.158

$ws.Cells.Item(123, 1).Value = -4
$ws.Cells.Item(123, 2).Value = 10
$ws.Cells.Item(123, 3).Value = 6.6888
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 250
$ws.Cells.Item(123, 6).Value = 112
$ws.Cells.Item(123, 7).Value = 0.16
$ws.Cells.Item(123, 8).Value = 0.158

$ws.Cells.Item(124, 1).Value = -3
$ws.Cells.Item(124, 2).Value = 10
$ws.Cells.Item(124, 3).Value = 6.8123
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 250
$ws.Cells.Item(124, 6).Value = 112
$ws.Cells.Item(124, 7).Value = 0.16
$ws.Cells.Item(124, 8).Value = 0.158

$ws.Cells.Item(125, 1).Value = 0
$ws.Cells.Item(125, 2).Value = 10
$ws.Cells.Item(125, 3).Value = 4.5922
$ws.Cells.Item(125, 4).Value = 0
$ws.Cells.Item(125, 5).Value = 250
$ws.Cells.Item(125, 6).Value = 112
$ws.Cells.Item(125, 7).Value = 0.16
$ws.Cells.Item(125, 8).Value = 0.158

$ws.Cells.Item(126, 1).Value = 1
$ws.Cells.Item(126, 2).Value = 10
$ws.Cells.Item(126, 3).Value = 4.6839
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 250
$ws.Cells.Item(126, 6).Value = 112
$ws.Cells.Item(126, 7).Value = 0.16
$ws.Cells.Item(126, 8).Value = 0.158

$ws.Cells.Item(127, 1).Value = 2
$ws.Cells.Item(127, 2).Value = 10
$ws.Cells.Item(127, 3).Value = 3.0886
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 5).Value = 250
$ws.Cells.Item(127, 6).Value = 112
$ws.Cells.Item(127, 7).Value = 0.16
$ws.Cells.Item(127, 8).Value = 0.158
